$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 6384
$ws.Range("B4").Value = 2050
$ws.Range("C4").Value = 11466
$ws.Range("D4").Value = 0.6917284967365771
$ws.Range("E4").Value = 0.2221613724737022
$ws.Range("F4").Value = 1.242381000740646
$ws.Range("G4").Value = 84423
$ws.Range("H4").Value = 35673
$ws.Range("I4").Value = 106604
$ws.Range("J4").Value = 0.02831751231846258
$ws.Range("K4").Value = 0.01522156619740772
$ws.Range("L4").Value = 0.03352731658531885
$ws.Range("M4").Value = 4647
$ws.Range("N4").Value = 1907
$ws.Range("O4").Value = 15549
$ws.Range("P4").Value = 0.5264324771852076
$ws.Range("Q4").Value = 0.2160054079508113
$ws.Range("R4").Value = 1.761334148472777
$ws.Range("S4").Value = 68556
$ws.Range("T4").Value = 33681
$ws.Range("U4").Value = 111189
$ws.Range("V4").Value = 0.03620315578577338
$ws.Range("W4").Value = 0.01839130680364554
$ws.Range("X4").Value = 0.05923561406248443
